# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 holds the "K" (strikeout) stat per start. The sheet had
# been populated with pitch-level "Strike#" counts; this regenerates that
# column with the real strikeout totals ("K") for each start, keeping every
# other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2  = 4
    3  = 3
    4  = 3
    5  = 2
    6  = 2
    7  = 4
    8  = 2
    9  = 4
    10 = 4
    11 = 3
    12 = 2
    13 = 6
    14 = 2
    15 = 2
    16 = 4
    17 = 5
    18 = 4
    19 = 3
    20 = 6
    21 = 6
    22 = 11
    23 = 6
    24 = 7
    25 = 4
    26 = 3
    27 = 6
    28 = 1
    29 = 4
    30 = 5
    31 = 4
    32 = 6
    33 = 9
    34 = 5
    35 = 3
    36 = 1
    37 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
